$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New column CB holds the "29-sep" data, one column after the existing last
# column CA ("28-sep").
$ws.Range("CB1").Value = "29-sep"

$values = @{
    "CB2"  = 12
    "CB3"  = 18
    "CB4"  = 10
    "CB5"  = 11
    "CB6"  = 9
    "CB7"  = 16
    "CB8"  = 12
    "CB9"  = 14
    "CB10" = 30
    "CB11" = 11
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

# Match formatting of the existing data columns (e.g. CA) for the new CB
# column: header uses the same style as CA1, data cells use the same style
# as CA2:CA11.
$ws.Range("CA1").Copy()
$ws.Range("CB1").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("CA2:CA11").Copy()
$ws.Range("CB2:CB11").PasteSpecial(-4122) # xlPasteFormats

$excel.CutCopyMode = 0

$ws.Range("CB12").Select()
